$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete column E ("align"); old column F ("zhuyin") shifts left to become E.
$ws.Columns("E:E").Delete()

# 2. Column B currently holds duplicate numeric row indices (leftover from the old
#    B column). Replace every B1:B62 cell with the "align" string that used to live
#    in column E (left / right / center / centerTitle), matching the new header "align".
$bValues = @(
    'align',
    'centerTitle',
    'left',
    'right',
    'left',
    'right',
    'center',
    'left',
    'right',
    'center',
    'left',
    'right',
    'center',
    'left',
    'right',
    'left',
    'right',
    'left',
    'right',
    'left',
    'right',
    'left',
    'right',
    'left',
    'right',
    'left',
    'right',
    'left',
    'right',
    'left',
    'right',
    'left',
    'right',
    'left',
    'right',
    'left',
    'right',
    'left',
    'right',
    'left',
    'right',
    'left',
    'right',
    'left',
    'right',
    'left',
    'right',
    'center',
    'left',
    'right',
    'left',
    'right',
    'left',
    'right',
    'left',
    'right',
    'left',
    'right',
    'left',
    'right',
    'left',
    'right'
)

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $bValues[$i]
}

# 3. Column widths: best-fit the five remaining columns (A:E) to their content.
$ws.Columns.Item(1).ColumnWidth = 3.5924479166666665
$ws.Columns.Item(2).ColumnWidth = 9.877604166666666
$ws.Columns.Item(3).ColumnWidth = 17.592447916666668
$ws.Columns.Item(4).ColumnWidth = 26.451822916666668
$ws.Columns.Item(5).ColumnWidth = 47.022135416666664

# 4. Restore the active-cell selection that was captured with the workbook.
$ws.Range("M13").Select()

